# Apply updated cryptocurrency price/volume data to Sheet1 (columns D and E).
# D-column numeric-looking values are entered with a leading apostrophe so
# Excel keeps them as text (matching the source data, which stores prices as
# text strings like '216.55' / '1.879.19' rather than numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.790.68"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "'1.649.02"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("E4").Value = "  +0.66%  "

$ws.Range("D5").Value = "'216.55"
$ws.Range("E5").Value = "  +1.31%  "

$ws.Range("E6").Value = "  +0.54%  "

$ws.Range("E7").Value = "  +0.56%  "

$ws.Range("E8").Value = "  +0.59%  "

$ws.Range("D10").Value = "'19.26"
$ws.Range("E10").Value = "  +1.86%  "

$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "'1.879.19"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").Value = "'1.651.15"
$ws.Range("E13").Value = "  +1.51%  "

$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("D16").Value = "'65.44"
$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("D17").Value = "'26.797.31"
$ws.Range("E17").Value = "  +0.84%  "

$ws.Range("D18").Value = "'0.0₃0745"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").Value = "'218.47"

$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").Value = "'2.52"
$ws.Range("E21").Value = "  +11.89%  "

$ws.Range("D22").Value = "'4.38"
$ws.Range("E22").Value = "  +1.34%  "

$ws.Range("E23").Value = "  +0.73%  "

$ws.Range("E24").Value = "  +1.21%  "

$ws.Range("D25").Value = "'145.94"
$ws.Range("E25").Value = "  -0.64%  "

$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "'7.21"
$ws.Range("E28").Value = "  +4.35%  "

$ws.Range("D29").Value = "'15.84"
$ws.Range("E29").Value = "  +1.09%  "

$ws.Range("D30").Value = "'0.0520"
$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("E31").Value = "  +1.42%  "

$ws.Range("D32").Value = "'3.36"
$ws.Range("E32").Value = "  +0.15%  "

$ws.Range("D33").Value = "'3.02"
$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("D34").Value = "'1.283.10"
$ws.Range("E34").Value = "  +1.16%  "

$ws.Range("E35").Value = "  +2.36%  "

$ws.Range("E37").Value = "  +1.90%  "

$ws.Range("D38").Value = "'0.542"
$ws.Range("E38").Value = "  +5.94%  "

$ws.Range("D39").Value = "'0.833"
$ws.Range("E39").Value = "  +4.06%  "

$ws.Range("E40").Value = "  +0.63%  "

$ws.Range("D41").Value = "'0.818"
$ws.Range("E41").Value = "  +2.51%  "

$ws.Range("E42").Value = "  -1.01%  "

$ws.Range("D43").Value = "'5.46"
$ws.Range("E43").Value = "  +1.84%  "

$ws.Range("D44").Value = "'1.790.47"
$ws.Range("E44").Value = "  +0.97%  "

$ws.Range("D45").Value = "'92.08"
$ws.Range("E45").Value = "  -1.51%  "

$ws.Range("D46").Value = "'59.89"
$ws.Range("E46").Value = "  +8.60%  "

$ws.Range("D47").Value = "'1.61"
$ws.Range("E47").Value = "  +0.66%  "

$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("E49").Value = "  +0.70%  "

$ws.Range("D50").Value = "'7.80"
$ws.Range("E50").Value = "  +2.58%  "

$ws.Range("D51").Value = "'0.0980"
$ws.Range("E51").Value = "  +1.82%  "
